# Corrected first year calendar error
# - Rows 2-41: Exam (H) flag corrected from 0 to 1
# - Rows 18-21: Holiday (F) corrected from 3 to 1
# - Rows 22-29: Holiday (F) corrected from 3 to 0
# - Rows 50-113: Holiday (F) corrected to 17
# - Rows 50-125: SemType (E) corrected from "Acad" to "Vacation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exam column (H) fix for rows 2-41
$ws.Range("H2:H41").Value = 1.0

# Holiday column (F) fixes
$ws.Range("F18:F21").Value = 1.0
$ws.Range("F22:F29").Value = 0.0
$ws.Range("F50:F113").Value = 17.0

# SemType column (E) fix: rows 50-125 were mis-marked "Acad", should be "Vacation"
$ws.Range("E50:E125").Value = "Vacation"
